$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''70.073.03'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '''3.505.04'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''605.63'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '''172.53'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('D8').Value = '''3.499.99'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '''0.194'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D11').Value = '''7.23'
$ws.Range('E11').Value = '  +7.05%  '
$ws.Range('D12').Value = '''0.587'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '''46.08'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').Value = '''4.078.79'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '''8.37'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '''612.41'
$ws.Range('E17').Value = '  -2.28%  '
$ws.Range('D18').Value = '''3.515.18'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '''70.162.23'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '''17.55'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').Value = '''0.878'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').Value = '''9.18'
$ws.Range('E23').Value = '  -8.14%  '
$ws.Range('D24').Value = '''98.63'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('D25').Value = '''15.54'
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('E26').Value = '  -3.52%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('D29').Value = '''33.80'
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('D30').Value = '''9.00'
$ws.Range('E30').Value = '  -2.55%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''8.06'
$ws.Range('E31').Value = '  -4.86%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').Value = '''2.98'
$ws.Range('E32').Value = '  -4.26%  '
$ws.Range('E33').Value = '  -4.93%  '
$ws.Range('D34').Value = '''631.91'
$ws.Range('E34').Value = '  +11.06%  '
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('D37').Value = '''10.74'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').Value = '''0.0479'
$ws.Range('E38').Value = '  +5.11%  '
$ws.Range('D39').Value = '''3.47'
$ws.Range('E39').Value = '  -5.19%  '
$ws.Range('D40').Value = '''56.77'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.146'
$ws.Range('E41').Value = '  +1.75%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '''3.361.98'
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('D44').Value = '''0.0₃0730'
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('E45').Value = '  -5.78%  '
$ws.Range('E46').Value = '  -4.65%  '
$ws.Range('D47').Value = '''31.88'
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('E48').Value = '  -4.16%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').Value = '''133.13'
$ws.Range('E50').Value = '  -0.57%  '
